$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 44.855544
$ws.Range("H2").Value = 134.566632
$ws.Range("I2").Value = 0.9269627513664965
$ws.Range("J2").Value = 0.9269627513664968
$ws.Range("M2").Value = 133.7780026666667
$ws.Range("N2").Value = 401.334008
$ws.Range("O2").Value = 0.50863533211804
$ws.Range("P2").Value = 0.5086353321180399
$ws.Range("Q2").Value = 6000.685084846785
$ws.Range("R2").Value = 54006.16576362106
$ws.Range("S2").Value = 0.4714860069023501
$ws.Range("T2").Value = 0.4714860069023502
$ws.Range("G3").Value = 44.855544
$ws.Range("H3").Value = 134.566632
$ws.Range("I3").Value = 0.9269627513664965
$ws.Range("J3").Value = 0.9269627513664968
$ws.Range("O3").Value = 0.1993888292903622
$ws.Range("P3").Value = 0.1993888292903622
$ws.Range("Q3").Value = 2352.313137637224
$ws.Range("R3").Value = 21170.81823873502
$ws.Range("S3").Value = 0.1848260177907389
$ws.Range("T3").Value = 0.1848260177907389
$ws.Range("G4").Value = 44.855544
$ws.Range("H4").Value = 134.566632
$ws.Range("I4").Value = 0.9269627513664965
$ws.Range("J4").Value = 0.9269627513664968
$ws.Range("M4").Value = 21.197691
$ws.Range("N4").Value = 63.593073
$ws.Range("O4").Value = 0.08059542216956049
$ws.Range("P4").Value = 0.08059542216956046
$ws.Range("Q4").Value = 950.8339613489042
$ws.Range("R4").Value = 8557.505652140137
$ws.Range("S4").Value = 0.07470895428184013
$ws.Range("T4").Value = 0.07470895428184011
$ws.Range("G5").Value = 44.855544
$ws.Range("H5").Value = 134.566632
$ws.Range("I5").Value = 0.9269627513664965
$ws.Range("J5").Value = 0.9269627513664968
$ws.Range("M5").Value = 55.59592133333333
$ws.Range("N5").Value = 166.787764
$ws.Range("O5").Value = 0.2113804164220374
$ws.Range("P5").Value = 0.2113804164220373
$ws.Range("Q5").Value = 2493.785295587872
$ws.Range("R5").Value = 22444.06766029085
$ws.Range("S5").Value = 0.1959417723915675
$ws.Range("T5").Value = 0.1959417723915675
$ws.Range("I6").Value = 0.03026428998407557
$ws.Range("J6").Value = 0.03026428998407558
$ws.Range("M6").Value = 133.7780026666667
$ws.Range("N6").Value = 401.334008
$ws.Range("O6").Value = 0.50863533211804
$ws.Range("P6").Value = 0.5086353321180399
$ws.Range("Q6").Value = 195.9156106792881
$ws.Range("R6").Value = 1763.240496113592
$ws.Range("S6").Value = 0.01539348718736695
$ws.Range("T6").Value = 0.01539348718736695
$ws.Range("I7").Value = 0.03026428998407557
$ws.Range("J7").Value = 0.03026428998407558
$ws.Range("O7").Value = 0.1993888292903622
$ws.Range("P7").Value = 0.1993888292903622
$ws.Range("Q7").Value = 76.80037501599301
$ws.Range("R7").Value = 691.2033751439371
$ws.Range("S7").Value = 0.006034361349228863
$ws.Range("T7").Value = 0.006034361349228864
$ws.Range("I8").Value = 0.03026428998407557
$ws.Range("J8").Value = 0.03026428998407558
$ws.Range("M8").Value = 21.197691
$ws.Range("N8").Value = 63.593073
$ws.Range("O8").Value = 0.08059542216956049
$ws.Range("P8").Value = 0.08059542216956046
$ws.Range("Q8").Value = 31.04365810875301
$ws.Range("R8").Value = 279.3929229787771
$ws.Range("S8").Value = 0.002439163227928572
$ws.Range("T8").Value = 0.002439163227928572
$ws.Range("I9").Value = 0.03026428998407557
$ws.Range("J9").Value = 0.03026428998407558
$ws.Range("M9").Value = 55.59592133333333
$ws.Range("N9").Value = 166.787764
$ws.Range("O9").Value = 0.2113804164220374
$ws.Range("P9").Value = 0.2113804164220373
$ws.Range("Q9").Value = 81.41928166200401
$ws.Range("R9").Value = 732.773534958036
$ws.Range("S9").Value = 0.006397278219551189
$ws.Range("T9").Value = 0.006397278219551188
$ws.Range("G10").Value = 1.967437666666666
$ws.Range("H10").Value = 5.902312999999999
$ws.Range("I10").Value = 0.04065810533109158
$ws.Range("J10").Value = 0.0406581053310916
$ws.Range("M10").Value = 133.7780026666667
$ws.Range("N10").Value = 401.334008
$ws.Range("O10").Value = 0.50863533211804
$ws.Range("P10").Value = 0.5086353321180399
$ws.Range("Q10").Value = 263.1998814178338
$ws.Range("R10").Value = 2368.798932760504
$ws.Range("S10").Value = 0.02068014890837002
$ws.Range("T10").Value = 0.02068014890837002
$ws.Range("G11").Value = 1.967437666666666
$ws.Range("H11").Value = 5.902312999999999
$ws.Range("I11").Value = 0.04065810533109158
$ws.Range("J11").Value = 0.0406581053310916
$ws.Range("O11").Value = 0.1993888292903622
$ws.Range("P11").Value = 0.1993888292903622
$ws.Range("Q11").Value = 103.176309059641
$ws.Range("R11").Value = 928.586781536769
$ws.Range("S11").Value = 0.008106772023130586
$ws.Range("T11").Value = 0.008106772023130588
$ws.Range("G12").Value = 1.967437666666666
$ws.Range("H12").Value = 5.902312999999999
$ws.Range("I12").Value = 0.04065810533109158
$ws.Range("J12").Value = 0.0406581053310916
$ws.Range("M12").Value = 21.197691
$ws.Range("N12").Value = 63.593073
$ws.Range("O12").Value = 0.08059542216956049
$ws.Range("P12").Value = 0.08059542216956046
$ws.Range("Q12").Value = 41.705135719761
$ws.Range("R12").Value = 375.346221477849
$ws.Range("S12").Value = 0.003276857163773784
$ws.Range("T12").Value = 0.003276857163773784
$ws.Range("G13").Value = 1.967437666666666
$ws.Range("H13").Value = 5.902312999999999
$ws.Range("I13").Value = 0.04065810533109158
$ws.Range("J13").Value = 0.0406581053310916
$ws.Range("M13").Value = 55.59592133333333
$ws.Range("N13").Value = 166.787764
$ws.Range("O13").Value = 0.2113804164220374
$ws.Range("P13").Value = 0.2113804164220373
$ws.Range("Q13").Value = 109.3815097442369
$ws.Range("R13").Value = 984.4335876981318
$ws.Range("S13").Value = 0.008594327235817196
$ws.Range("T13").Value = 0.008594327235817198
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.1023373333333333
$ws.Range("H14").Value = 0.307012
$ws.Range("I14").Value = 0.002114853318336234
$ws.Range("J14").Value = 0.002114853318336234
$ws.Range("M14").Value = 133.7780026666667
$ws.Range("N14").Value = 401.334008
$ws.Range("O14").Value = 0.50863533211804
$ws.Range("P14").Value = 0.5086353321180399
$ws.Range("Q14").Value = 13.69048405156622
$ws.Range("R14").Value = 123.214356464096
$ws.Range("S14").Value = 0.001075689119952889
$ws.Range("T14").Value = 0.001075689119952889
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.1023373333333333
$ws.Range("H15").Value = 0.307012
$ws.Range("I15").Value = 0.002114853318336234
$ws.Range("J15").Value = 0.002114853318336234
$ws.Range("O15").Value = 0.1993888292903622
$ws.Range("P15").Value = 0.1993888292903622
$ws.Range("Q15").Value = 5.366771466884001
$ws.Range("R15").Value = 48.300943201956
$ws.Range("S15").Value = 0.0004216781272638994
$ws.Range("T15").Value = 0.0004216781272638995
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.1023373333333333
$ws.Range("H16").Value = 0.307012
$ws.Range("I16").Value = 0.002114853318336234
$ws.Range("J16").Value = 0.002114853318336234
$ws.Range("M16").Value = 21.197691
$ws.Range("N16").Value = 63.593073
$ws.Range("O16").Value = 0.08059542216956049
$ws.Range("P16").Value = 0.08059542216956046
$ws.Range("Q16").Value = 2.169315169764
$ws.Range("R16").Value = 19.523836527876
$ws.Range("S16").Value = 0.0001704474960180047
$ws.Range("T16").Value = 0.0001704474960180047
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.1023373333333333
$ws.Range("H17").Value = 0.307012
$ws.Range("I17").Value = 0.002114853318336234
$ws.Range("J17").Value = 0.002114853318336234
$ws.Range("M17").Value = 55.59592133333333
$ws.Range("N17").Value = 166.787764
$ws.Range("O17").Value = 0.2113804164220374
$ws.Range("P17").Value = 0.2113804164220373
$ws.Range("Q17").Value = 5.689538333463111
$ws.Range("R17").Value = 51.20584500116799
$ws.Range("S17").Value = 0.0004470385751014407
$ws.Range("T17").Value = 0.0004470385751014407
